# Card11: add a new service event (row 18) and fill in the previously
# blank cells of the prior event (row 17) with the literal text "nan"
# (matching how the rest of the sheet represents missing values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card11")

# --- Row 17: the existing last row had several untyped/blank cells;
# the commit fills them with the literal string "nan".
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","M","P")) {
    $ws.Range($col + "17").Value = "nan"
}

# --- Row 18: brand-new event row appended below.
# Column A holds the card number "11" stored as text (it already is
# text throughout the sheet), so force text formatting before writing
# it to avoid Excel auto-converting it to a number.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "11"

# Columns B-K, M and P are left blank (empty text) for this event, same
# as most other rows in the sheet.
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","M","P")) {
    $ws.Range($col + "18").Value = "'"
}

# Date / correction / serviced-by details for the new event.
$ws.Range("L18").Value = "19\10\2024"
$ws.Range("N18").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O18").Value = "تيم العمل"
